$wb = $excel.ActiveWorkbook

# ---------- Sheet 1: Statistics ----------
$ws1 = $wb.Worksheets.Item("Statistics")

$stats = @(
    ,@('2', '2024-08-21 17:43:46', 53.2597868721498, 7)
    ,@('3', '2024-08-21 17:43:48', 53.3096261348934, 13)
    ,@('4', '2024-08-21 17:43:50', 52.88633371650054, 17)
    ,@('5', '2024-08-21 17:43:52', 50.84237885885849, 20)
    ,@('6', '2024-08-21 17:43:55', 45.85710907513953, 26)
    ,@('7', '2024-08-21 17:43:57', 39.65579930384778, 27)
    ,@('8', '2024-08-21 17:43:59', 40.37960121286115, 28)
    ,@('9', '2024-08-21 17:44:01', 38.29742760463628, 27)
    ,@('10', '2024-08-21 17:44:03', 36.75830900142348, 29)
    ,@('11', '2024-08-21 17:44:05', 39.2270736932723, 29)
    ,@('12', '2024-08-21 17:44:07', 34.96951146285748, 30)
    ,@('13', '2024-08-21 17:44:09', 31.93928113526318, 31)
    ,@('14', '2024-08-21 17:44:11', 30.87915370635895, 31)
    ,@('15', '2024-08-21 17:44:13', 30.86046802726679, 30)
    ,@('16', '2024-08-21 17:44:15', 34.8145223948564, 31)
    ,@('17', '2024-08-21 17:44:17', 32.47670200564293, 31)
    ,@('18', '2024-08-21 17:44:19', 33.66226437536362, 31)
    ,@('19', '2024-08-21 17:44:21', 33.09571951800906, 34)
    ,@('20', '2024-08-21 17:44:23', 38.59116503744637, 31)
    ,@('21', '2024-08-21 17:44:25', 34.69362189409538, 30)
    ,@('22', '2024-08-21 17:44:27', 34.10056881538849, 27)
    ,@('23', '2024-08-21 17:44:29', 33.66405266894405, 28)
    ,@('24', '2024-08-21 17:44:31', 26.17422947859876, 29)
    ,@('25', '2024-08-21 17:44:33', 32.31625505677929, 28)
    ,@('26', '2024-08-21 17:44:35', 27.34599918758866, 30)
    ,@('27', '2024-08-21 17:44:37', 25.20376652458257, 29)
    ,@('28', '2024-08-21 17:44:39', 34.65633244540913, 30)
    ,@('29', '2024-08-21 17:44:41', 32.02488864699885, 29)
    ,@('30', '2024-08-21 17:44:43', 31.7561179230039, 27)
    ,@('31', '2024-08-21 17:44:45', 34.44263519006502, 28)
    ,@('32', '2024-08-21 17:44:47', 28.91155524530594, 33)
    ,@('33', '2024-08-21 17:44:49', 31.2247647270954, 27)
    ,@('34', '2024-08-21 17:44:51', 33.31848549298379, 29)
    ,@('35', '2024-08-21 17:44:53', 33.28308946787771, 31)
    ,@('36', '2024-08-21 17:44:55', 31.51743266704381, 30)
    ,@('37', '2024-08-21 17:44:57', 27.54575450675285, 31)
    ,@('38', '2024-08-21 17:44:59', 29.62615157164711, 29)
    ,@('39', '2024-08-21 17:45:01', 32.13586844328807, 32)
    ,@('40', '2024-08-21 17:45:03', 31.95387765513365, 32)
    ,@('41', '2024-08-21 17:45:05', 34.9106350205066, 31)
    ,@('42', '2024-08-21 17:45:07', 30.52806368336728, 34)
    ,@('43', '2024-08-21 17:45:09', 25.09747607378401, 33)
    ,@('44', '2024-08-21 17:45:11', 25.77450295835182, 32)
    ,@('45', '2024-08-21 17:45:13', 27.00591701907633, 33)
    ,@('46', '2024-08-21 17:45:15', 28.41774230506413, 34)
    ,@('47', '2024-08-21 17:45:17', 27.2860076742514, 35)
    ,@('48', '2024-08-21 17:45:19', 29.31387494205724, 38)
    ,@('49', '2024-08-21 17:45:21', 26.64486446228091, 37)
    ,@('50', '2024-08-21 17:45:23', 26.73714009450719, 39)
    ,@('51', '2024-08-21 17:45:25', 24.564870681794, 41)
    ,@('52', '2024-08-21 17:45:27', 25.10075636602006, 42)
    ,@('53', '2024-08-21 17:45:29', 22.29508633083137, 42)
    ,@('54', '2024-08-21 17:45:31', 19.38186955297007, 41)
    ,@('55', '2024-08-21 17:45:33', 17.6735899755845, 42)
    ,@('56', '2024-08-21 17:45:35', 18.46172362378111, 44)
    ,@('57', '2024-08-21 17:45:37', 16.83440958380139, 43)
    ,@('58', '2024-08-21 17:45:40', 15.72603118014481, 44)
)

foreach ($r in $stats) {
    $rowNum = $r[0]
    $ws1.Range("A$rowNum").Value = $r[1]
    $ws1.Range("B$rowNum").Value = [double]$r[2]
    $ws1.Range("C$rowNum").Value = [double]$r[3]
}

# ---------- Sheet 2: Accidents ----------
$ws2 = $wb.Worksheets.Item("Accidents")

$acc = @(
    ,@('2', '2024-08-21 17:43:54', 'Car and Car', '45.36 and 56.25', 1)
    ,@('3', '2024-08-21 17:43:57', 'Car and Car', '0.00 and 52.72', 2)
    ,@('4', '2024-08-21 17:43:59', 'Car and Car', '51.88 and 40.09', 3)
    ,@('5', '2024-08-21 17:44:01', 'Car and Truck', '28.73 and 47.09', 4)
    ,@('6', '2024-08-21 17:44:05', 'Car and Car', '42.41 and 47.40', 5)
    ,@('7', '2024-08-21 17:44:06', 'Car and Car', '0.00 and 40.57', 6)
    ,@('8', '2024-08-21 17:44:09', 'Car and Car', '3.03 and 48.20', 7)
    ,@('9', '2024-08-21 17:44:10', 'Car and Car', '34.15 and 61.49', 8)
    ,@('10', '2024-08-21 17:44:14', 'Car and Car', '57.60 and 51.34', 9)
    ,@('11', '2024-08-21 17:44:15', 'Truck and Car', '36.20 and 50.52', 10)
    ,@('12', '2024-08-21 17:44:18', 'Car and Car', '0.00 and 23.76', 11)
    ,@('13', '2024-08-21 17:44:18', 'Car and Car', '65.22 and 40.38', 12)
    ,@('14', '2024-08-21 17:44:20', 'Car and Car', '32.49 and 0.00', 11)
    ,@('15', '2024-08-21 17:44:21', 'Car and Car', '0.00 and 58.69', 13)
    ,@('16', '2024-08-21 17:44:23', 'Car and Car', '48.41 and 37.20', 14)
    ,@('17', '2024-08-21 17:44:29', 'Car and Car', '55.84 and 45.41', 15)
    ,@('18', '2024-08-21 17:44:29', 'Car and Car', '0.00 and 59.09', 16)
    ,@('19', '2024-08-21 17:44:34', 'Car and Car', '43.81 and 44.49', 17)
    ,@('20', '2024-08-21 17:44:34', 'Car and Car', '33.39 and 62.28', 18)
    ,@('21', '2024-08-21 17:44:39', 'Car and Car', '38.69 and 53.66', 19)
    ,@('22', '2024-08-21 17:44:40', 'Car and Car', '57.29 and 52.54', 20)
    ,@('23', '2024-08-21 17:44:45', 'Car and Car', '63.71 and 44.16', 21)
    ,@('24', '2024-08-21 17:44:45', 'Car and Truck', '52.51 and 31.50', 22)
    ,@('25', '2024-08-21 17:44:50', 'Car and Car', '54.36 and 45.10', 23)
    ,@('26', '2024-08-21 17:44:54', 'Truck and Car', '0.00 and 46.19', 24)
    ,@('27', '2024-08-21 17:44:56', 'Car and Car', '47.22 and 46.27', 25)
    ,@('28', '2024-08-21 17:44:59', 'Car and Truck', '55.15 and 40.70', 26)
    ,@('29', '2024-08-21 17:45:01', 'Car and Car', '56.16 and 47.61', 27)
    ,@('30', '2024-08-21 17:45:07', 'Truck and Car', '38.43 and 55.51', 28)
)

foreach ($r in $acc) {
    $rowNum = $r[0]
    $ws2.Range("A$rowNum").Value = $r[1]
    $ws2.Range("B$rowNum").Value = $r[2]
    $ws2.Range("C$rowNum").Value = $r[3]
    $ws2.Range("D$rowNum").Value = [double]$r[4]
}
